$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 36 - ORTEGA ROMAN KLEBER ERWIN : PANELES DECORATIVOS (K)
$wsVentasGrupo.Range("K36").Value = 263.16

# Row 37 - ORTEGA ROMAN LUIS FERNANDO : LAVABOS (I), PIEDRA SINTERIZADA (L)
$wsVentasGrupo.Range("I37").Value = 156.6
$wsVentasGrupo.Range("L37").Value = 556.8099999999999

# Row 56 - counts of advisors meeting target ("X de 54")
$wsVentasGrupo.Range("I56").Value = "5 de 54"
$wsVentasGrupo.Range("K56").Value = "1 de 54"
$wsVentasGrupo.Range("L56").Value = "8 de 54"

# --- Sheet "VENTA MENSUAL" ----------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 36 - ORTEGA ROMAN KLEBER ERWIN : julio (F)
$wsVentaMensual.Range("F36").Value = 14683.88

# Row 37 - ORTEGA ROMAN LUIS FERNANDO : julio (F)
$wsVentaMensual.Range("F37").Value = 3813.7

# Row 56 - totals : julio (F)
$wsVentaMensual.Range("F56").Value = 76307.53

# --- Sheet "CUMPLIMIENTO MENSUAL" --------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 8 - LAVABOS
$wsCumplimiento.Range("D8").Value = 666.76
$wsCumplimiento.Range("E8").Value = 333.24
$wsCumplimiento.Range("F8").Value = 0.66676

# Row 12 - PANELES DECORATIVOS
$wsCumplimiento.Range("D12").Value = 263.16
$wsCumplimiento.Range("E12").Value = 86.83999999999997
$wsCumplimiento.Range("F12").Value = 0.7518857142857144

# Row 15 - PIEDRA SINTERIZADA
$wsCumplimiento.Range("D15").Value = 6869.49
$wsCumplimiento.Range("E15").Value = 6630.51
$wsCumplimiento.Range("F15").Value = 0.5088511111111111

# Row 19 - TOTAL
$wsCumplimiento.Range("D19").Value = 76307.53
$wsCumplimiento.Range("E19").Value = 37398.92064517915
$wsCumplimiento.Range("F19").Value = 0.6710923572675536
